# ncp-gop-transect-2020-info.xlsx edit
# "edit projects, temperature defn, KM's role"
#
# 1. ColumnHeadersNcp: fix the "temperature" attribute's definition text,
#    which erroneously described salinity instead of temperature.
# 2. Personnel: add a new row for Kate Morkeski (metadata Provider).

$wb = $excel.ActiveWorkbook

# --- 1. ColumnHeadersNcp: correct the temperature attribute definition ---
$wsNcp = $wb.Worksheets.Item("ColumnHeadersNcp")
$wsNcp.Range("B7").Value = "Underway thermosalinograph temperature in degrees Celsius. URI http://vocab.nerc.ac.uk/collection/P01/current/TEMPSZ01/"
$wsNcp.Activate()
$wsNcp.Range("B7").Select()

# --- 2. Personnel: add Kate Morkeski as metadata Provider ---
$wsPersonnel = $wb.Worksheets.Item("Personnel")
$wsPersonnel.Range("A10").Value = "Kate"
$wsPersonnel.Range("C10").Value = "Morkeski"
$wsPersonnel.Range("D10").Value = "Northeast U.S. Shelf LTER"
$wsPersonnel.Range("E10").Value = "kmorkeski@whoi.edu"
$wsPersonnel.Range("F10").Value = "0000-0002-2903-5851"
$wsPersonnel.Range("G10").Value = "metadata Provider"
$wsPersonnel.Range("H10").Value = "Northeast U.S. Shelf LTER"
$wsPersonnel.Range("I10").Value = "NSF"
$wsPersonnel.Range("J10").Value = "OCE-2322676"
$wsPersonnel.Activate()
$wsPersonnel.Range("A10:J10").Select()
